$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's header row + data row (originally rows 1-2) need to move down
# by one row, freeing row 1 (used for "alternate receptors" bookkeeping).
# Inserting a blank row above row 1 shifts the existing rows 1-2 down to 2-3,
# carrying their values/styles along automatically.
$ws.Rows.Item(1).Insert()

# Reflect the resulting active selection (matches the saved file's cursor).
$ws.Range("C8").Select()
